$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# Fill in the newly-tracked weekly burn-down data for columns F (week 6) and G (week 7)
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0

$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0

$ws.Range("F6").Value = 10
$ws.Range("G6").Value = 0

$ws.Range("F7").Value = 10
$ws.Range("G7").Value = 0

$ws.Range("F8").Value = 10
$ws.Range("G8").Value = 10

$ws.Range("F9").Value = 10
$ws.Range("G9").Value = 10

$ws.Range("F10").Value = 10
$ws.Range("G10").Value = 10

$ws.Range("F11").Value = 10
$ws.Range("G11").Value = 10

# Overwrite the old shared formula in E13 with its computed literal value
$ws.Range("E13").Value = 70

# Populate the actual burn-down total for the newly-tracked week
$ws.Range("G13").Value = 40

# Move the selection/active cell as it ended up after the edits, and reset the scrolled view
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("F14").Select()
